$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Tanggal Bayar" header column (I) -----------------------------
# Set the new header text first so its shared-string entry is allocated
# before the other content edits below (keeps shared-string order stable).
$ws.Range("I1").Value = "Tanggal Bayar"

# Match the header formatting already used by the other header cells
# (bold font, no special color) instead of whatever default style a bare
# Value assignment would pick up.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Tanggal Bayar"

# --- Updated payer / payment content -----------------------------------
$ws.Range("B2").Value = "Anastasia Novitasari"
$ws.Range("B3").Value = "Anastasia Novitasari"
$ws.Range("E2").Value = "SPP1"
$ws.Range("E3").Value = "SPP1"
$ws.Range("F2").Value = "Januari"
$ws.Range("F3").Value = "Februari"
$ws.Range("G3").Value = 50000

# --- New "Tanggal Bayar" values (import/transaction dates) --------------
$ws.Range("I2").NumberFormat = "mm-dd-yy"
$ws.Range("I2").Value = (Get-Date -Year 2021 -Month 6 -Day 12 -Hour 9 -Minute 31 -Second 0)

# Reuse I2's exact cell format for I3 (copy/paste-format) so both date
# cells share a single style entry instead of minting a duplicate one.
$ws.Range("I2").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = (Get-Date -Year 2021 -Month 6 -Day 13 -Hour 15 -Minute 20 -Second 0)

# --- Column width for the new column ------------------------------------
$ws.Columns.Item(9).ColumnWidth = 12.26

# --- Selection -----------------------------------------------------------
$ws.Range("F3").Select()

$excel.CutCopyMode = $false
